# Update countries & provincias Spain
#
# Refreshes the COVID-19 country table with newer totals and re-inserts
# "Ghana" with updated figures (it previously sorted after "Oman" with
# stale numbers; with the refreshed total it now ranks above
# "Luxemburgo", pushing Luxemburgo/Afganistan/Nigeria/Hungria/Oman down
# one row each). "Belice" and "Nueva Caledonia" swap places as well.
# Also bumps the "last updated" timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados..." footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 8 de Mayo de 2020 a las 22:04"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1314295
$ws.Range("C4").Value = 21672
$ws.Range("D4").Value = 220997
$ws.Range("E4").Value = 1015050
$ws.Range("F4").Value = 16772
$ws.Range("G4").Value = 1320
$ws.Range("H4").Value = 78248

# Francia (row 9)
$ws.Range("B9").Value = 175462
$ws.Range("C9").Value = 671
$ws.Range("E9").Value = 93450

# Barein (row 61)
$ws.Range("B61").Value = 4444
$ws.Range("C61").Value = 245
$ws.Range("D61").Value = 2028
$ws.Range("E61").Value = 2408

# Ghana moves up into row 62 (ahead of Luxemburgo) with refreshed data
$ws.Range("A62").Value = "Ghana"
$ws.Range("B62").Value = 4012
$ws.Range("C62").Value = 921
$ws.Range("D62").Value = 323
$ws.Range("E62").Value = 3671
$ws.Range("F62").Value = 4
$ws.Range("H62").Value = 18

# Luxemburgo shifts down to row 63 (values unchanged)
$ws.Range("A63").Value = "Luxemburgo"
$ws.Range("B63").Value = 3871
$ws.Range("C63").Value = 12
$ws.Range("D63").Value = 3526
$ws.Range("E63").Value = 245
$ws.Range("F63").Value = 16
$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 100

# Afganistan shifts down to row 64 (values unchanged)
$ws.Range("A64").Value = "Afganistan"
$ws.Range("B64").Value = 3778
$ws.Range("C64").Value = 215
$ws.Range("D64").Value = 472
$ws.Range("E64").Value = 3197
$ws.Range("F64").Value = 7
$ws.Range("G64").Value = 3
$ws.Range("H64").Value = 109

# Nigeria shifts down to row 65 (values unchanged)
$ws.Range("A65").Value = "Nigeria"
$ws.Range("B65").Value = 3526
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 601
$ws.Range("E65").Value = 2818
$ws.Range("F65").Value = 4
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 107

# Hungria shifts down to row 66 (values unchanged)
$ws.Range("A66").Value = "Hungria"
$ws.Range("B66").Value = 3178
$ws.Range("C66").Value = 28
$ws.Range("D66").Value = 865
$ws.Range("E66").Value = 1921
$ws.Range("F66").Value = 74
$ws.Range("G66").Value = 9
$ws.Range("H66").Value = 392

# Oman shifts down to row 67 (values unchanged)
$ws.Range("A67").Value = "Oman"
$ws.Range("B67").Value = 3112
$ws.Range("C67").Value = 154
$ws.Range("D67").Value = 1025
$ws.Range("E67").Value = 2071
$ws.Range("F67").Value = 17
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 16

# Nepal (row 158)
$ws.Range("B158").Value = 102
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 31
$ws.Range("E158").Value = 71

# Belice now sorts ahead of Nueva Caledonia (row 192)
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Nueva Caledonia shifts down to row 193
$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0
